$wb = $excel.ActiveWorkbook

# --- Update selection on the "Produkt" sheet (D22 -> G2) ---
$produkt = $wb.Worksheets.Item("Produkt")
[void]$produkt.Range("G2").Select()

# --- Add the new "Variablen" worksheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$variablen = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$variablen.Name = "Variablen"

# --- Populate the new sheet with the Name/Wert key-value pair ---
$variablen.Range("A1").Value = "Name"
$variablen.Range("A2").Value = "Wert"
$variablen.Range("B1").Value = "Rücksendekosten"

# "0.1" needs to be stored as text (shared string), not as a number.
$variablen.Range("B2").NumberFormat = "@"
$variablen.Range("B2").Value = "0.1"
$variablen.Range("B2").Style = "Normal"

# --- Column widths to roughly match the original layout ---
$variablen.Columns.Item(1).ColumnWidth = 14.83
$variablen.Columns.Item(2).ColumnWidth = 14.33

# --- Final selection on the new sheet ---
[void]$variablen.Range("B2").Select()
